$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.294417977333069
$ws.Range("B1").Value = 2.045036554336548
$ws.Range("C1").Value = 5.383988857269287
$ws.Range("D1").Value = 1.91789710521698
$ws.Range("E1").Value = 1.096792459487915
